$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ts1 = "2025-10-17T07:09:31.429905"
$ts2 = "2025-10-17T07:09:31.498985"
$ts3 = "2025-10-17T07:09:31.569376"
$ts4 = "2025-10-17T07:09:31.570374"
$ts5 = "2025-10-17T07:09:31.571374"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 26).Value = $ts1
}
for ($r = 16; $r -le 25; $r++) {
    $ws.Cells.Item($r, 26).Value = $ts2
}
for ($r = 26; $r -le 32; $r++) {
    $ws.Cells.Item($r, 26).Value = $ts3
}
for ($r = 33; $r -le 41; $r++) {
    $ws.Cells.Item($r, 26).Value = $ts4
}
for ($r = 42; $r -le 48; $r++) {
    $ws.Cells.Item($r, 26).Value = $ts5
}
